$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 30 (this shifts rows 30-37 down
# to 31-38, preserving all their existing data/formatting).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price entry.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44988
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = 100112003
$ws.Range("G30").Value = "Ajo"
$ws.Range("H30").Value = "Chino"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 700
$ws.Range("K30").Value = 16000
$ws.Range("L30").Value = 17000
$ws.Range("M30").Value = 16357
$ws.Range("N30").Value = "$/caja 10 kilos"
$ws.Range("O30").Value = "China"
$ws.Range("P30").Value = 1636
$ws.Range("Q30").Value = 10
$ws.Range("R30").Value = "Hortaliza"
